# Apply the changes described by the diff:
# 1. Rename sheet "op2" -> "wong3"
# 2. Update row 4 values: B4, C4, E4, F4 from 1 -> 2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet
$ws.Name = "wong3"

# 2. Update the data cells in row 4
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2
